$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("P1").Value = "Oier Talavera"
